# Automatische test-sync: 2025-06-17 22:08:15
# Append a new "Afmelding nieuwsbrief" (newsletter unsubscribe) entry to the
# mail log on the "Logs" sheet, and bump the matching tally on "Dashboard".

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# New row of log data (row 46).
$newRow = 46
$logs.Cells.Item($newRow, 1).Value = "Afmelding nieuwsbrief"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$logs.Cells.Item($newRow, 4).Value = "Afmelding"
$logs.Cells.Item($newRow, 6).Value = "2025-06-17 22:07:25"
$logs.Cells.Item($newRow, 7).Value = "Nee"

# Extend the conditional formatting ranges to cover the new row, mirroring
# how Excel keeps the Categorie/Beantwoord highlighting in sync.
$categoryConditions = $logs.Range("D2:D45").FormatConditions
foreach ($cond in $categoryConditions) {
    $cond.ModifyAppliesToRange($logs.Range("D2:D46"))
}

$repliedConditions = $logs.Range("G2:G45").FormatConditions
foreach ($cond in $repliedConditions) {
    $cond.ModifyAppliesToRange($logs.Range("G2:G46"))
}

# Update the Dashboard summary count for the "Afmelding" category (was 7).
$dashboard.Cells.Item(4, 2).Value = 8
